$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Footer "datetimeFigureOut" date fields: "1/24/20" -> "3/23/21"
#    (handout master, notes master, slide master, and all 11 slide layouts)
# ---------------------------------------------------------------------------
function Update-DateShape($shape) {
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq "1/24/20") {
            $tr.Text = "3/23/21"
        }
    }
}

# Handout master
$hm = $p.HandoutMaster
for ($j = 1; $j -le $hm.Shapes.Count; $j++) {
    Update-DateShape $hm.Shapes.Item($j)
}

# Notes master
$nm = $p.NotesMaster
for ($j = 1; $j -le $nm.Shapes.Count; $j++) {
    Update-DateShape $nm.Shapes.Item($j)
}

# Slide master
$m = $p.SlideMaster
for ($j = 1; $j -le $m.Shapes.Count; $j++) {
    Update-DateShape $m.Shapes.Item($j)
}

# All slide layouts
$layouts = $m.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $lyt = $layouts.Item($i)
    for ($j = 1; $j -le $lyt.Shapes.Count; $j++) {
        Update-DateShape $lyt.Shapes.Item($j)
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 2 ("Course contents"): swap order of the first two bullet lines
#    "A formula for problem-solving simple algorithm problems" / "How to read code"
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(2)
$tr2 = $sh2.TextFrame2.TextRange
$para1 = $tr2.Paragraphs(1, 1)
$para2 = $tr2.Paragraphs(2, 1)
$para1.Text = "How to read code"
$para2.Text = "A formula for problem-solving simple algorithm problems"

# ---------------------------------------------------------------------------
# 3) Slide 3 ("Course projects"): percentage + grading text updates
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(2)
$tr3 = $sh3.TextFrame.TextRange

# 3a) "kmeans++ initial point selection (20%)" -> "...(17%)"
$full = $tr3.Text
$needle = "++ initial point selection (20%)"
$idx = $full.IndexOf($needle)
$rng = $tr3.Characters($idx + 1, $needle.Length)
$rng.Text = "++ initial point selection (17%)"

# 3b) "Feature importance and selection (22%)" -> "...(20%)"
$full = $tr3.Text
$needle = "Feature importance and selection (22%)"
$idx = $full.IndexOf($needle)
$rng = $tr3.Characters($idx + 1, $needle.Length)
$rng.Text = "Feature importance and selection (20%)"

# 3c) "grader will assign check -, check, check+ based upon your reports"
#     -> two runs: "I " + "will assign check, check-, check-- based upon your reports"
$full = $tr3.Text
$needle = "grader will assign check -, check, check+ based upon your reports"
$idx = $full.IndexOf($needle)
$rng = $tr3.Characters($idx + 1, $needle.Length)
$rng.Text = "I will assign check, check-, check-- based upon your reports"

$full = $tr3.Text
$needle2 = "I will assign check, check-, check-- based upon your reports"
$idx2 = $full.IndexOf($needle2)
$rngSplit = $tr3.Characters($idx2 + 1, 2)
$rngSplit.Text = "I "
